$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.242.68"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "1.998.67"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").Value = "0.7775"
$ws.Range("E5").Value = "  +29.96%  "
$ws.Range("D6").Value = "256.54"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "0.3472"
$ws.Range("E8").Value = "  +13.92%  "
$ws.Range("D9").Value = "28.14"
$ws.Range("E9").Value = "  +15.69%  "
$ws.Range("D10").Value = "0.07266"
$ws.Range("E10").Value = "  +9.20%  "
$ws.Range("D11").Value = "0.8449"
$ws.Range("E11").Value = "  +6.46%  "
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "100.91"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.000.37"
$ws.Range("E14").Value = "  +4.51%  "
$ws.Range("D15").Value = "5.647"
$ws.Range("E15").Value = "  +6.04%  "
$ws.Range("D16").Value = "15.69"
$ws.Range("E16").Value = "  +16.04%  "
$ws.Range("D17").Value = "272.55"
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "31.267.95"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "'0.000008336"
$ws.Range("E19").Value = "  +9.46%  "
$ws.Range("D20").Value = "6.007"
$ws.Range("E20").Value = "  +9.96%  "
$ws.Range("D21").Value = "2.255.52"
$ws.Range("E21").Value = "  +5.54%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "7.131"
$ws.Range("E24").Value = "  +8.64%  "
$ws.Range("E25").Value = "  +8.82%  "
$ws.Range("D26").Value = "164.73"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "0.1413"
$ws.Range("E27").Value = "  +31.80%  "
$ws.Range("D28").Value = "'20.00"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").Value = "'2.410"
$ws.Range("E29").Value = "  +23.01%  "
$ws.Range("D30").Value = "'1.610"
$ws.Range("E30").Value = "  +5.75%  "
$ws.Range("D31").Value = "'4.660"
$ws.Range("E31").Value = "  +5.84%  "
$ws.Range("D32").Value = "1.368"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "4.484"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").Value = "0.05375"
$ws.Range("E34").Value = "  +9.01%  "
$ws.Range("D35").Value = "1.265"
$ws.Range("E35").Value = "  +9.52%  "
$ws.Range("D36").Value = "0.7891"
$ws.Range("E36").Value = "  +10.53%  "
$ws.Range("D37").Value = "'2.780"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "'0.9990"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").Value = "0.02015"
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").Value = "2.942"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "85.78"
$ws.Range("E41").Value = "  +11.48%  "
$ws.Range("D42").Value = "6.842"
$ws.Range("E42").Value = "  +7.28%  "
$ws.Range("D43").Value = "0.4701"
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("D44").Value = "2.139"
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("D45").Value = "'0.8630"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "'105.30"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("D47").Value = "10.23"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "7.776"
$ws.Range("E49").Value = "  +7.99%  "
$ws.Range("D50").Value = "'37.90"
$ws.Range("E50").Value = "  +6.58%  "
$ws.Range("D51").Value = "'3.000"
$ws.Range("E51").Value = "  +42.33%  "
